$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 2).Value = "NONE"
    $ws.Cells.Item($r, 4).Value = "NONE"
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4)).Style = "Normal"
}
